$d = $word.ActiveDocument

# 1. Merge the four runs "o.s.v." / " " / "Hemsidan ... kommer " / "minska ... kvaliteten. "
#    into a single run with the concatenated text (same visible text, one run).
$mergedText = "o.s.v. Hemsidan är också väldigt stor p.g.a. att det finns så många filmer, vilket betyder att jag kommer minska antalet sidor som jag designar för att i stället fokusera på kvaliteten. "
$rng = $d.Content
$null = $rng.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)

# 2. Left-align the three milestone-table paragraphs that currently inherit the
#    table's centered alignment: "Mappstruktur och grundstrukturen av html klar",
#    "Inlagt många av de element ..." and "Design av hemsidan".
#    (Range.Paragraphs items returned from Find are not stable handles, so walk
#    Document.Paragraphs directly and match on the paragraph's own text.)
$targets = @(
    "Mappstruktur och grundstrukturen av html klar",
    "Inlagt många av de element som hemsidan ska bestå av.",
    "Design av hemsidan"
)
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -ne $null) {
        foreach ($target in $targets) {
            if ($t.StartsWith($target)) {
                $p.Alignment = 0
            }
        }
    }
}
